$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projects")

$ws.Range("A2").Value = "testeObregon"
$ws.Range("A3").Value = "AmxCoMovPosPlan090"
$ws.Range("A4").ClearContents()

$ws.Range("A3").Select()
